$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 889
$ws1.Cells.Item(5, 6).Value = 1189
$ws1.Cells.Item(6, 6).Value = 69
$ws1.Cells.Item(7, 6).Value = 4364
$ws1.Cells.Item(8, 6).Value = 2600
$ws1.Cells.Item(10, 6).Value = 2514
$ws1.Cells.Item(14, 6).Value = 1660
$ws1.Cells.Item(15, 6).Value = 662
$ws1.Cells.Item(16, 6).Value = 35
$ws1.Cells.Item(18, 6).Value = 325
$ws1.Cells.Item(22, 6).Value = 25
$ws1.Cells.Item(23, 6).Value = 480
$ws1.Cells.Item(26, 6).Value = 547
$ws1.Cells.Item(27, 6).Value = 691
$ws1.Cells.Item(28, 6).Value = 105
$ws1.Cells.Item(33, 6).Value = 1026
$ws1.Cells.Item(34, 6).Value = 124
$ws1.Cells.Item(36, 6).Value = 1127
$ws1.Cells.Item(37, 6).Value = 2048
$ws1.Cells.Item(38, 6).Value = 265
$ws1.Cells.Item(39, 6).Value = 10
$ws1.Cells.Item(40, 6).Value = 544
$ws1.Cells.Item(42, 6).Value = 24
$ws1.Cells.Item(43, 6).Value = 658
$ws1.Cells.Item(44, 6).Value = 1325
$ws1.Cells.Item(45, 6).Value = 94
$ws1.Cells.Item(47, 6).Value = 432

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(5, 6).Value = 68

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 890
$ws4.Cells.Item(3, 6).Value = 1189
$ws4.Cells.Item(5, 6).Value = 69
$ws4.Cells.Item(6, 6).Value = 4364
$ws4.Cells.Item(7, 6).Value = 2600
$ws4.Cells.Item(8, 6).Value = 2514
$ws4.Cells.Item(9, 6).Value = 1660
$ws4.Cells.Item(12, 6).Value = 662
$ws4.Cells.Item(13, 6).Value = 35
$ws4.Cells.Item(15, 6).Value = 325
$ws4.Cells.Item(19, 6).Value = 480
$ws4.Cells.Item(22, 6).Value = 547
$ws4.Cells.Item(23, 6).Value = 691
$ws4.Cells.Item(24, 6).Value = 105
$ws4.Cells.Item(25, 6).Value = 68
$ws4.Cells.Item(31, 6).Value = 1026
$ws4.Cells.Item(32, 6).Value = 124
$ws4.Cells.Item(35, 6).Value = 2048
$ws4.Cells.Item(36, 6).Value = 265
$ws4.Cells.Item(40, 6).Value = 544
$ws4.Cells.Item(42, 6).Value = 24
$ws4.Cells.Item(43, 6).Value = 658
$ws4.Cells.Item(44, 6).Value = 1325
$ws4.Cells.Item(46, 6).Value = 94
$ws4.Cells.Item(47, 6).Value = 432
